$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Bid Amount cell (D2) to hold a text value "1,00,001" instead of the numeric 10000
$ws.Range("D2").Value = "1,00,001"

# Move the active selection from D5 to D4
$ws.Range("D4").Select()
